$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 83: EARNED credit of 1.25 for Nov 2022 period (date already present)
$ws1.Range("C83").Value = 1.25

# Row 84: Dec 2022 period
$ws1.Range("A84").Value = 44896
$ws1.Range("C84").Value = 1.25

# Row 85: new "2023" year header row, matching the style of other year rows
# (bold run inside an otherwise date-formatted / quote-prefixed cell)
$ws1.Range("A85").Value = "'2023"
$ws1.Range("A85").Characters(1,3).Font.Bold = $true
$ws1.Range("A85").Characters(4,1).Font.Bold = $true

# Row 86: Jan 2023 period
$ws1.Range("A86").Value = 44927
$ws1.Range("C86").Value = 1.25

# Row 87: Feb 2023 period
$ws1.Range("A87").Value = 44958
$ws1.Range("C87").Value = 1.25

# Row 88: Mar 2023 period
$ws1.Range("A88").Value = 44986
$ws1.Range("C88").Value = 1.25

# Row 89: Apr 2023 period
$ws1.Range("A89").Value = 45017
$ws1.Range("C89").Value = 1.25

# Row 90: May 2023 period
$ws1.Range("A90").Value = 45047
$ws1.Range("C90").Value = 1.25

# Row 91: Jun 2023 period
$ws1.Range("A91").Value = 45078
$ws1.Range("C91").Value = 1.25

# Row 92: Jul 2023 period, with an absence/undertime debit and an end-of-period date stamp
$ws1.Range("A92").Value = 45108
$ws1.Range("C92").Value = 1.25
$ws1.Range("H92").Value = 1
$ws1.Range("K21").Copy()
$ws1.Range("K92").PasteSpecial(-4122)
$ws1.Range("K92").Value = 45138

# Row 93: Aug 2023 period (date only)
$ws1.Range("A93").Value = 45139

# Row 94: Sep 2023 period (date only)
$ws1.Range("A94").Value = 45170
